$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Column D (Price) updates - forced to text to avoid numeric auto-conversion
Set-TextValue $ws.Range('D2') '61.203.75'
Set-TextValue $ws.Range('D3') '2.932.65'
Set-TextValue $ws.Range('D5') '593.03'
Set-TextValue $ws.Range('D6') '144.94'
Set-TextValue $ws.Range('D8') '0.504'
Set-TextValue $ws.Range('D9') '7.01'
Set-TextValue $ws.Range('D13') '33.71'
Set-TextValue $ws.Range('D15') '3.419.50'
Set-TextValue $ws.Range('D16') '61.186.61'
Set-TextValue $ws.Range('D18') '2.935.11'
Set-TextValue $ws.Range('D19') '433.52'
Set-TextValue $ws.Range('D20') '13.52'
Set-TextValue $ws.Range('D24') '11.09'
Set-TextValue $ws.Range('D26') '11.88'
Set-TextValue $ws.Range('D29') '2.61'
Set-TextValue $ws.Range('D30') '6.99'
Set-TextValue $ws.Range('D31') '0.110'
Set-TextValue $ws.Range('D32') '26.78'
Set-TextValue $ws.Range('D34') '0.0₃0878'
Set-TextValue $ws.Range('D37') '2.98'
Set-TextValue $ws.Range('D40') '8.63'
Set-TextValue $ws.Range('D41') '41.44'
Set-TextValue $ws.Range('D42') '0.283'
Set-TextValue $ws.Range('D43') '376.39'
Set-TextValue $ws.Range('D45') '2.703.47'
Set-TextValue $ws.Range('D46') '133.28'
Set-TextValue $ws.Range('D48') '23.94'

# Column E (Volume 1h) updates - plain text assignment is safe
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +4.13%  '
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('E15').Value = '  +0.87%  '
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -3.55%  '
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  -1.17%  '
$ws.Range('E31').Value = '  +3.15%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  +2.87%  '
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('E41').Value = '  +3.29%  '
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('E51').Value = '  +0.31%  '
